$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "to sing (～を)"
$ws.Range("A3").Value = "to put on (a hat) (～を)"
$ws.Range("A4").Value = "to put on (items below your waist) (～を)"
$ws.Range("A5").Value = "to get to know (～を)"
$ws.Range("A8").Value = "to live (～に)"
$ws.Range("A13").Value = "to put on (clothes above your waist) (～を)"
$ws.Range("A16").Value = "to get married (～と)"
